$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cell, $text) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-CellText "D2" "329.19"
Set-CellText "E2" "0.15%"
Set-CellText "D3" "44.28"
Set-CellText "E3" "0.90%"
Set-CellText "D4" "5.587"
Set-CellText "E4" "3.00%"
Set-CellText "D5" "0.08086"
Set-CellText "E5" "-0.24%"
Set-CellText "D6" "1.992"
Set-CellText "E6" "4.84%"
Set-CellText "D7" "0.9530"
Set-CellText "E7" "0.99%"
Set-CellText "E8" "-7.37%"
Set-CellText "D9" "0.1172"
Set-CellText "E9" "-1.36%"
Set-CellText "D10" "0.1853"
Set-CellText "E10" "-1.87%"
Set-CellText "B11" "MCDex"
Set-CellText "C11" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-CellText "D11" "10.18"
Set-CellText "E11" "18.35%"
Set-CellText "B12" "MandalaExchangeToken"
Set-CellText "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-CellText "D12" "0.09788"
Set-CellText "E12" "1.74%"
Set-CellText "B13" "BitrueCoin"
Set-CellText "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-CellText "D13" "0.04540"
Set-CellText "E13" "7.59%"
Set-CellText "B14" "BitMartToken"
Set-CellText "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-CellText "D14" "0.1068"
Set-CellText "E14" "-0.22%"
Set-CellText "B15" "BitForexToken"
Set-CellText "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-CellText "D15" "0.001281"
Set-CellText "E15" "0.44%"
Set-CellText "B16" "CoinExToken"
Set-CellText "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-CellText "D16" "0.04192"
Set-CellText "E16" "-4.53%"
Set-CellText "B17" "TigerCash"
Set-CellText "C17" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-CellText "D17" "0.005871"
Set-CellText "E17" "-1.63%"
Set-CellText "D18" "3.378"
Set-CellText "E18" "-5.21%"
Set-CellText "D19" "4.311"
Set-CellText "E19" "0.17%"
Set-CellText "D20" "0.3481"
Set-CellText "E20" "-1.02%"
Set-CellText "B21" "ProBitToken"
Set-CellText "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-CellText "D21" "0.1420"
Set-CellText "E21" "4.47%"
Set-CellText "B22" "ZBToken"
Set-CellText "C22" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-CellText "D22" "0.2505"
Set-CellText "E22" "-3.81%"
Set-CellText "B23" "BitKan"
Set-CellText "C23" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-CellText "D23" "0.001244"
Set-CellText "E23" "0.23%"
Set-CellText "B24" "HotbitToken"
Set-CellText "C24" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-CellText "D24" "0.004356"
Set-CellText "E24" "0.65%"
Set-CellText "D25" "0.0001189"
Set-CellText "E26" "-0.94%"
Set-CellText "D38" "0.02684"
Set-CellText "E38" "-0.63%"
Set-CellText "D39" "0.05570"
Set-CellText "E39" "0.47%"
Set-CellText "D40" "0.007574"
Set-CellText "E40" "-2.94%"
Set-CellText "D41" "0.1410"
Set-CellText "E41" "0.87%"
Set-CellText "D42" "0.007948"
Set-CellText "E42" "-18.59%"
Set-CellText "D43" "0.002014"
Set-CellText "E43" "-5.30%"
Set-CellText "D44" "0.008408"
Set-CellText "E44" "-12.76%"
Set-CellText "D45" "0.00007189"
Set-CellText "E45" "1.13%"
Set-CellText "D46" "0.00000000749"
Set-CellText "E46" "-0.76%"
Set-CellText "B47" "BOLO"
Set-CellText "C47" "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-CellText "D47" "0.004387"
Set-CellText "E47" "26.19%"
Set-CellText "B48" "CoinbaseStockToken"
Set-CellText "C48" "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-CellText "D48" "0.002269"
Set-CellText "E48" "-0.75%"
Set-CellText "D49" "0.00002099"
Set-CellText "E49" "-0.76%"
Set-CellText "D50" "0.0001999"
Set-CellText "E50" "-0.76%"
